$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "815×4="
$t.Cell(1,2).Range.Text = "174×9="
$t.Cell(1,3).Range.Text = "413×8="
$t.Cell(1,4).Range.Text = "625×9="
$t.Cell(1,5).Range.Text = "706×9="
$t.Cell(5,1).Range.Text = "571×5="
$t.Cell(5,2).Range.Text = "521×2="
$t.Cell(5,3).Range.Text = "739×2="
$t.Cell(5,4).Range.Text = "741×9="
$t.Cell(5,5).Range.Text = "880×3="
$t.Cell(10,1).Range.Text = "285×3="
$t.Cell(10,2).Range.Text = "835×8="
$t.Cell(10,3).Range.Text = "211×8="
$t.Cell(10,4).Range.Text = "327×5="
$t.Cell(10,5).Range.Text = "597×7="
$t.Cell(15,1).Range.Text = "599×7="
$t.Cell(15,2).Range.Text = "159×5="
$t.Cell(15,3).Range.Text = "531×6="
$t.Cell(15,4).Range.Text = "413×7="
$t.Cell(15,5).Range.Text = "743×5="
$t.Cell(20,1).Range.Text = "963×4="
$t.Cell(20,2).Range.Text = "673×7="
$t.Cell(20,3).Range.Text = "904×2="
$t.Cell(20,4).Range.Text = "133×8="
$t.Cell(20,5).Range.Text = "735×7="
